$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1851851851851852
$ws.Range("C2").Value = 0.5648148148148148
$ws.Range("J2").Value = 0.006172839506172839
$ws.Range("P2").Value = 0.1450617283950617
$ws.Range("S2").Value = 0.09876543209876543
$ws.Range("B3").Value = 0.01522842639593909
$ws.Range("C3").Value = 0.05583756345177665
$ws.Range("J3").Value = 0.03045685279187817
$ws.Range("P3").Value = 0.7157360406091371
$ws.Range("S3").Value = 0.182741116751269
$ws.Range("J4").Value = 0.103448275862069
$ws.Range("P4").Value = 0.4827586206896552
$ws.Range("S4").Value = 0.4137931034482759
$ws.Range("B6").Value = 0.06425702811244979
$ws.Range("D6").Value = 0.008032128514056224
$ws.Range("F6").Value = 0.06827309236947791
$ws.Range("J6").Value = 0.2971887550200803
$ws.Range("O6").Value = 0.02008032128514056
$ws.Range("Q6").Value = 0.1124497991967871
$ws.Range("R6").Value = 0.06024096385542169
$ws.Range("S6").Value = 0.3694779116465863
$ws.Range("B7").Value = 0.1042944785276074
$ws.Range("D7").Value = 0.0245398773006135
$ws.Range("F7").Value = 0.03680981595092025
$ws.Range("J7").Value = 0.1165644171779141
$ws.Range("O7").Value = 0.03680981595092025
$ws.Range("Q7").Value = 0.1533742331288344
$ws.Range("R7").Value = 0.1104294478527607
$ws.Range("S7").Value = 0.4171779141104294
$ws.Range("B8").Value = 0.09555555555555556
$ws.Range("D8").Value = 0.008888888888888889
$ws.Range("E8").Value = 0.002222222222222222
$ws.Range("F8").Value = 0.07111111111111111
$ws.Range("J8").Value = 0.1066666666666667
$ws.Range("O8").Value = 0.01333333333333333
$ws.Range("Q8").Value = 0.1888888888888889
$ws.Range("R8").Value = 0.09777777777777778
$ws.Range("S8").Value = 0.4155555555555556
$ws.Range("B9").Value = 0.1297071129707113
$ws.Range("D9").Value = 0.008368200836820083
$ws.Range("F9").Value = 0.08368200836820083
$ws.Range("J9").Value = 0.07949790794979079
$ws.Range("O9").Value = 0.02510460251046025
$ws.Range("Q9").Value = 0.1799163179916318
$ws.Range("R9").Value = 0.07949790794979079
$ws.Range("S9").Value = 0.4142259414225942
$ws.Range("B10").Value = 0.1277340332458443
$ws.Range("D10").Value = 0.01662292213473316
$ws.Range("E10").Value = 0.001749781277340333
$ws.Range("F10").Value = 0.07874015748031496
$ws.Range("J10").Value = 0.09798775153105861
$ws.Range("O10").Value = 0.02012248468941382
$ws.Range("Q10").Value = 0.1846019247594051
$ws.Range("R10").Value = 0.08573928258967629
$ws.Range("S10").Value = 0.3867016622922135
$ws.Range("G11").Value = 0.1553784860557769
$ws.Range("J11").Value = 0.07569721115537849
$ws.Range("K11").Value = 0.1832669322709163
$ws.Range("L11").Value = 0.5537848605577689
$ws.Range("S11").Value = 0.03187250996015936
$ws.Range("G12").Value = 0.7866666666666666
$ws.Range("J12").Value = 0.1733333333333333
$ws.Range("K12").Value = 0.006666666666666667
$ws.Range("L12").Value = 0.006666666666666667
$ws.Range("S12").Value = 0.02666666666666667
$ws.Range("G13").Value = 0.696969696969697
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("F15").Value = 0.05333333333333334
$ws.Range("H15").Value = 0.1111111111111111
$ws.Range("I15").Value = 0.08888888888888889
$ws.Range("J15").Value = 0.3777777777777778
$ws.Range("K15").Value = 0.05777777777777778
$ws.Range("O15").Value = 0.06222222222222222
$ws.Range("S15").Value = 0.2488888888888889
$ws.Range("F16").Value = 0.02604166666666667
$ws.Range("H16").Value = 0.2083333333333333
$ws.Range("I16").Value = 0.08854166666666667
$ws.Range("J16").Value = 0.3385416666666667
$ws.Range("K16").Value = 0.08854166666666667
$ws.Range("M16").Value = 0.01041666666666667
$ws.Range("N16").Value = 0.005208333333333333
$ws.Range("O16").Value = 0.07291666666666667
$ws.Range("S16").Value = 0.1614583333333333
$ws.Range("F17").Value = 0.02267002518891688
$ws.Range("H17").Value = 0.1486146095717884
$ws.Range("I17").Value = 0.1209068010075567
$ws.Range("J17").Value = 0.3979848866498741
$ws.Range("K17").Value = 0.0906801007556675
$ws.Range("M17").Value = 0.02015113350125945
$ws.Range("O17").Value = 0.05289672544080604
$ws.Range("S17").Value = 0.146095717884131
$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.1727748691099476
$ws.Range("I18").Value = 0.1204188481675393
$ws.Range("J18").Value = 0.4031413612565445
$ws.Range("K18").Value = 0.09424083769633508
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.03664921465968586
$ws.Range("S18").Value = 0.1413612565445026
$ws.Range("F19").Value = 0.02262090483619345
$ws.Range("H19").Value = 0.2285491419656786
$ws.Range("I19").Value = 0.1029641185647426
$ws.Range("J19").Value = 0.3424336973478939
$ws.Range("K19").Value = 0.08814352574102964
$ws.Range("M19").Value = 0.01638065522620905
$ws.Range("N19").Value = 0.0007800312012480499
$ws.Range("O19").Value = 0.07800312012480499
$ws.Range("S19").Value = 0.1201248049921997
